# Update the "check_ssl_cert" stats workbook after the latest release:
# append a new day's worth of statistics as row 25 of the "Data" table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")
$lo = $ws.ListObjects.Item(1)

# Add a new row to the "Data" table; this also grows the table range
# (and the worksheet dimension) from A1:X24 to A1:X25.
$newRow = $lo.ListRows.Add()
$rng = $newRow.Range

# Copy the formatting (number formats, styles, ...) from the previous
# last row so the new row renders identically to the existing rows.
$ws.Range("A24:X24").Copy()
$rng.PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the new row of stats.
$rng.Item(1, 1).Value = 44723    # Date
$rng.Item(1, 2).Value = 105      # Authors
$rng.Item(1, 3).Value = 236      # Versions
$rng.Item(1, 4).Value = 189      # GH Releases
$rng.Item(1, 5).Value = 4510     # LoC
$rng.Item(1, 6).Value = 1733     # Commits
$rng.Item(1, 7).Value = 3579     # File Changes
$rng.Item(1, 8).Value = 61141    # Insertions
$rng.Item(1, 9).Value = 43154    # Deletions
$rng.Item(1, 10).Value = 0       # Open issues
$rng.Item(1, 11).Value = 0       # Open bugs
$rng.Item(1, 12).Value = 221     # Closed issues
$rng.Item(1, 13).Formula = "=Data[[#This Row],[Open issues]]+Data[[#This Row],[Closed issues]]"   # Issues
$rng.Item(1, 14).Value = 0       # Open pull requests
$rng.Item(1, 15).Value = 158     # Closed pull requests
$rng.Item(1, 16).Formula = "=Data[[#This Row],[Open pull requests]]+Data[[#This Row],[Closed pull requests]]"   # Pull requests
$rng.Item(1, 17).Value = 131     # Command line options
$rng.Item(1, 18).Value = 132     # Tests
$rng.Item(1, 19).Value = 4       # GH workflows
$rng.Item(1, 20).Value = 0       # Running
$rng.Item(1, 21).Value = 304     # Failed
$rng.Item(1, 22).Value = 501     # OK
$rng.Item(1, 23).Value = 7       # Cancelled
$rng.Item(1, 24).Formula = "=SUM(Data[[#This Row],[Running]:[Cancelled]])"   # GH runs

# Reflect the author's final selection/scroll position on the Data sheet.
$ws.Activate()
$ws.Range("T26").Select()
